$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.008.82'
$ws.Range('E2').Value = '  -2.04%  '

$ws.Range('D3').Value = '''1.792.71'
$ws.Range('E3').Value = '  -2.70%  '

$ws.Range('D4').Value = '''1.011'
$ws.Range('E4').Value = '  +0.84%  '

$ws.Range('E5').Value = '  +0.43%  '

$ws.Range('D6').Value = '''307.36'
$ws.Range('E6').Value = '  -2.02%  '

$ws.Range('D7').Value = '''0.4153'
$ws.Range('E7').Value = '  -2.26%  '

$ws.Range('D8').Value = '''0.3533'
$ws.Range('E8').Value = '  -3.66%  '

$ws.Range('D9').Value = '''0.06998'
$ws.Range('E9').Value = '  -3.38%  '

$ws.Range('D10').Value = '''0.8382'
$ws.Range('E10').Value = '  -3.69%  '

$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '''1.922.27'
$ws.Range('E11').Value = '  -0.12%  '

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = '''20.03'
$ws.Range('E12').Value = '  -3.74%  '

$ws.Range('D13').Value = '''5.223'
$ws.Range('E13').Value = '  -3.18%  '

$ws.Range('D14').Value = '''6.300'
$ws.Range('E14').Value = '  -3.62%  '

$ws.Range('D15').Value = '''0.06810'
$ws.Range('E15').Value = '  -2.15%  '

$ws.Range('D16').Value = '''1.011'
$ws.Range('E16').Value = '  +0.75%  '

$ws.Range('D17').Value = '''79.38'
$ws.Range('E17').Value = '  -0.86%  '

$ws.Range('D18').Value = '''0.000008661'
$ws.Range('E18').Value = '  -3.87%  '

$ws.Range('E19').Value = '  +0.41%  '

$ws.Range('D20').Value = '''14.97'
$ws.Range('E20').Value = '  -3.34%  '

$ws.Range('D21').Value = '''27.413.57'
$ws.Range('E21').Value = '  -1.34%  '

$ws.Range('D22').Value = '''5.021'
$ws.Range('E22').Value = '  -0.65%  '

$ws.Range('D23').Value = '''10.63'
$ws.Range('E23').Value = '  -2.51%  '

$ws.Range('D24').Value = '''2.062.68'
$ws.Range('E24').Value = '  -2.51%  '

$ws.Range('D25').Value = '''1.946'
$ws.Range('E25').Value = '  -0.95%  '

$ws.Range('D26').Value = '''152.31'
$ws.Range('E26').Value = '  -1.25%  '

$ws.Range('D27').Value = '''18.04'
$ws.Range('E27').Value = '  -2.00%  '

$ws.Range('D28').Value = '''4.984'
$ws.Range('E28').Value = '  -5.42%  '

$ws.Range('D29').Value = '''112.10'
$ws.Range('E29').Value = '  -2.88%  '

$ws.Range('D30').Value = '''1.638'
$ws.Range('E30').Value = '  -10.85%  '

$ws.Range('D31').Value = '''0.08828'
$ws.Range('E31').Value = '  -0.61%  '

$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').Value = '''2.870'
$ws.Range('E32').Value = '  -2.90%  '

$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '''0.7150'
$ws.Range('E33').Value = '  -7.99%  '

$ws.Range('D34').Value = '''4.311'
$ws.Range('E34').Value = '  -5.70%  '

$ws.Range('D35').Value = '''1.006'
$ws.Range('E35').Value = '  +0.43%  '

$ws.Range('D36').Value = '''1.067'
$ws.Range('E36').Value = '  -7.73%  '

$ws.Range('D37').Value = '''1.070'
$ws.Range('E37').Value = '  -2.97%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.01880'
$ws.Range('E38').Value = '  -3.57%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.05063'
$ws.Range('E39').Value = '  -5.90%  '

$ws.Range('D40').Value = '''0.4899'
$ws.Range('E40').Value = '  -4.68%  '

$ws.Range('D41').Value = '''0.1602'
$ws.Range('E41').Value = '  -3.52%  '

$ws.Range('D42').Value = '''2.626'
$ws.Range('E42').Value = '  -6.97%  '

$ws.Range('D43').Value = '''6.131'
$ws.Range('E43').Value = '  -10.52%  '

$ws.Range('D44').Value = '''7.994'
$ws.Range('E44').Value = '  -6.07%  '

$ws.Range('D45').Value = '''1.006'
$ws.Range('E45').Value = '  +0.27%  '

$ws.Range('D46').Value = '''10.21'
$ws.Range('E46').Value = '  -2.85%  '

$ws.Range('D47').Value = '''103.51'
$ws.Range('E47').Value = '  -2.55%  '

$ws.Range('D48').Value = '''0.06308'
$ws.Range('E48').Value = '  -3.60%  '

$ws.Range('D49').Value = '''0.4499'
$ws.Range('E49').Value = '  -4.62%  '

$ws.Range('D50').Value = '''1.572'
$ws.Range('E50').Value = '  -3.66%  '

$ws.Range('D51').Value = '''61.79'
$ws.Range('E51').Value = '  -4.11%  '

Write-Host "Edit complete"
